$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 799, shifting the existing rows 799:840 down to 800:841.
$ws.Rows.Item(799).Insert()

# Populate the newly inserted row with the new data point.
# Column A holds a date-like string that must remain literal text (not be
# auto-converted to a date serial number by Excel's input parser), so we
# temporarily force a text number format, assign the value, then clear the
# formatting again so the cell ends up with no explicit style - matching
# the rest of the data rows in the sheet.
$ws.Range("A799").NumberFormat = "@"
$ws.Range("A799").Value = "2026/02/12"
$ws.Range("A799").ClearFormats()

$ws.Range("B799").Value = "木"
$ws.Range("C799").Value = 7
$ws.Range("D799").Value = 201
